$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Copy format (style) of the last existing data row (539) down to the new rows (540:553)
$ws.Range("A539:C539").Copy()
$ws.Range("A540:C553").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(540, 1).Value = 'cs'
$ws.Cells.Item(540, 2).Value = 'lab.build.preview.description'
$ws.Cells.Item(540, 3).Value = 'Popis'

$ws.Cells.Item(541, 1).Value = 'cs'
$ws.Cells.Item(541, 2).Value = 'lab.build.preview.atomizer'
$ws.Cells.Item(541, 3).Value = 'Atomizér'

$ws.Cells.Item(542, 1).Value = 'cs'
$ws.Cells.Item(542, 2).Value = 'lab.build.preview.coil'
$ws.Cells.Item(542, 3).Value = 'Spirálka'

$ws.Cells.Item(543, 1).Value = 'cs'
$ws.Cells.Item(543, 2).Value = 'lab.build.preview.cotton'
$ws.Cells.Item(543, 3).Value = 'Vata'

$ws.Cells.Item(544, 1).Value = 'cs'
$ws.Cells.Item(544, 2).Value = 'lab.build.preview.ohm'
$ws.Cells.Item(544, 3).Value = 'Odpor'

$ws.Cells.Item(545, 1).Value = 'cs'
$ws.Cells.Item(545, 2).Value = 'lab.build.preview.coilOffset'
$ws.Cells.Item(545, 3).Value = 'Pozice spirálky'

$ws.Cells.Item(546, 1).Value = 'cs'
$ws.Cells.Item(546, 2).Value = 'lab.build.preview.cottonOffset'
$ws.Cells.Item(546, 3).Value = 'Množství vaty'

$ws.Cells.Item(547, 1).Value = 'cs'
$ws.Cells.Item(547, 2).Value = 'lab.build.preview.coils'
$ws.Cells.Item(547, 3).Value = 'Počet spirálek'

$ws.Cells.Item(548, 1).Value = 'cs'
$ws.Cells.Item(548, 2).Value = 'lab.build.button.clone'
$ws.Cells.Item(548, 3).Value = 'Klonovat'

$ws.Cells.Item(549, 1).Value = 'cs'
$ws.Cells.Item(549, 2).Value = 'lab.build.button.index'
$ws.Cells.Item(549, 3).Value = 'Detail buildu'

$ws.Cells.Item(550, 1).Value = 'cs'
$ws.Cells.Item(550, 2).Value = 'lab.build.preview'
$ws.Cells.Item(550, 3).Value = 'Náhled buildu'

$ws.Cells.Item(551, 1).Value = 'cs'
$ws.Cells.Item(551, 2).Value = 'lab.build.preview.preview.title'
$ws.Cells.Item(551, 3).Value = 'Náhled buildu'

$ws.Cells.Item(552, 1).Value = 'cs'
$ws.Cells.Item(552, 2).Value = 'lab.build.preview.preview.subtitle'
$ws.Cells.Item(552, 3).Value = 'Zde vidíte veškeré dostupné informace o buildu.'

$ws.Cells.Item(553, 1).Value = 'cs'
$ws.Cells.Item(553, 2).Value = 'lab.build.clone.title'
$ws.Cells.Item(553, 3).Value = 'Klon buildu'

# Update the active sheet selection / scroll position to reflect the new last edited cell
$ws.Activate()
$ws.Range("B548").Select()

